$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tests")

# The row "S'assurer que l'utilisateur garde le contrôle du système lors des
# interactions" (row 15) was removed from the test sheet. Deleting the entire
# row shifts every following row up by one, shrinks the used dimension/table
# range accordingly, and drops the row's text from the shared-string table
# when Excel re-saves the workbook.
$ws.Rows.Item(15).Delete()

# Reflect the cursor / selection position left behind after the edit.
$ws.Activate()
$ws.Range("A29").Select()
